# Re-run of the Katalon UM-Data test suite on Fri Jul 12 2024.
# Updates the "Date" column (B) timestamps for every test-run row
# across the workbook, and populates three newly-passed runs
# (rows 3-5) on CreateUserSCFNameErr / CreateUserSCLNameErr with
# Result="Pass" + a new Date.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("CreateUser")
$ws.Range("B2").Value = "Fri Jul 12 18:51:44 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 18:52:20 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 18:52:52 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserPasswordSpChar")
$ws.Range("B2").Value = "Fri Jul 12 19:03:26 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 19:03:55 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 19:04:25 EDT 2024"
$ws.Range("B5").Value = "Fri Jul 12 19:04:54 EDT 2024"
$ws.Range("B6").Value = "Fri Jul 12 19:05:24 EDT 2024"
$ws.Range("B7").Value = "Fri Jul 12 19:05:54 EDT 2024"
$ws.Range("B8").Value = "Fri Jul 12 19:06:24 EDT 2024"
$ws.Range("B9").Value = "Fri Jul 12 19:06:54 EDT 2024"
$ws.Range("B10").Value = "Fri Jul 12 19:07:24 EDT 2024"
$ws.Range("B11").Value = "Fri Jul 12 19:07:54 EDT 2024"
$ws.Range("B12").Value = "Fri Jul 12 19:08:23 EDT 2024"
$ws.Range("B13").Value = "Fri Jul 12 19:08:53 EDT 2024"
$ws.Range("B14").Value = "Fri Jul 12 19:09:22 EDT 2024"
$ws.Range("B15").Value = "Fri Jul 12 19:09:53 EDT 2024"
$ws.Range("B16").Value = "Fri Jul 12 19:10:22 EDT 2024"
$ws.Range("B17").Value = "Fri Jul 12 19:10:54 EDT 2024"

$ws = $wb.Worksheets.Item("ModifyUser")
$ws.Range("B2").Value = "Fri Jul 12 19:16:28 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 19:17:08 EDT 2024"

$ws = $wb.Worksheets.Item("ModifyUserPwd")
$ws.Range("B2").Value = "Fri Jul 12 19:12:49 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 19:13:21 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 19:13:53 EDT 2024"
$ws.Range("B5").Value = "Fri Jul 12 19:14:25 EDT 2024"
$ws.Range("B6").Value = "Fri Jul 12 19:14:57 EDT 2024"
$ws.Range("B7").Value = "Fri Jul 12 19:15:28 EDT 2024"
$ws.Range("B8").Value = "Fri Jul 12 19:15:55 EDT 2024"

$ws = $wb.Worksheets.Item("AddDeleteRole")
$ws.Range("B2").Value = "Fri Jul 12 19:22:46 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 19:23:14 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 19:23:40 EDT 2024"
$ws.Range("B5").Value = "Fri Jul 12 19:24:06 EDT 2024"

$ws = $wb.Worksheets.Item("SearchRole")
$ws.Range("B2").Value = "Fri Jul 12 19:24:32 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 19:24:57 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 19:25:19 EDT 2024"
$ws.Range("B5").Value = "Fri Jul 12 19:25:41 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserSpCharError")
$ws.Range("B2").Value = "Fri Jul 12 18:54:56 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 18:55:19 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 18:55:42 EDT 2024"
$ws.Range("B5").Value = "Fri Jul 12 18:56:04 EDT 2024"
$ws.Range("B6").Value = "Fri Jul 12 18:56:27 EDT 2024"
$ws.Range("B7").Value = "Fri Jul 12 18:56:50 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserSCFNameErr")
$ws.Range("B2").Value = "Fri Jul 12 19:17:49 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserSCLNameErr")
$ws.Range("B2").Value = "Fri Jul 12 19:19:21 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserSpChar")
$ws.Range("B2").Value = "Fri Jul 12 18:53:24 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 18:53:54 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 18:54:25 EDT 2024"

$ws = $wb.Worksheets.Item("UsernameCase")
$ws.Range("B2").Value = "Fri Jul 12 19:02:40 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 19:02:55 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 19:03:10 EDT 2024"

$ws = $wb.Worksheets.Item("PassCase")
$ws.Range("B2").Value = "Fri Jul 12 19:02:31 EDT 2024"

$ws = $wb.Worksheets.Item("FindUser")
$ws.Range("B2").Value = "Fri Jul 12 18:57:13 EDT 2024"
$ws.Range("B3").Value = "Fri Jul 12 18:57:43 EDT 2024"
$ws.Range("B4").Value = "Fri Jul 12 18:58:09 EDT 2024"
$ws.Range("B5").Value = "Fri Jul 12 18:58:35 EDT 2024"
$ws.Range("B6").Value = "Fri Jul 12 18:59:03 EDT 2024"
$ws.Range("B7").Value = "Fri Jul 12 18:59:27 EDT 2024"
$ws.Range("B8").Value = "Fri Jul 12 18:59:53 EDT 2024"
$ws.Range("B9").Value = "Fri Jul 12 19:00:19 EDT 2024"
$ws.Range("B10").Value = "Fri Jul 12 19:00:46 EDT 2024"
$ws.Range("B11").Value = "Fri Jul 12 19:01:14 EDT 2024"
$ws.Range("B12").Value = "Fri Jul 12 19:01:39 EDT 2024"
$ws.Range("B13").Value = "Fri Jul 12 19:02:05 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserErrors")
$ws.Range("B13").Value = "Fri Jul 12 19:20:54 EDT 2024"
$ws.Range("B14").Value = "Fri Jul 12 19:21:16 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserSCFNameErr")
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Jul 12 19:18:14 EDT 2024"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Jul 12 19:18:36 EDT 2024"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Jul 12 19:18:58 EDT 2024"

$ws = $wb.Worksheets.Item("CreateUserSCLNameErr")
$ws.Range("A3").Value = "Pass"
$ws.Range("B3").Value = "Fri Jul 12 19:19:45 EDT 2024"
$ws.Range("A4").Value = "Pass"
$ws.Range("B4").Value = "Fri Jul 12 19:20:07 EDT 2024"
$ws.Range("A5").Value = "Pass"
$ws.Range("B5").Value = "Fri Jul 12 19:20:30 EDT 2024"
